# Auto-generated edit script updating Leve market-price / profit figures
# across all 8 class sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Values correspond to refreshed market data pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 846.6042
$ws.Range("J17").Value = 844.6667
$ws.Range("L17").Value = 2534.0001
$ws.Range("N17").Value = -2870.0001

$ws.Range("H132").Value = 148481.44
$ws.Range("I132").Value = 157713.53
$ws.Range("K132").Value = 473140.59
$ws.Range("M132").Value = -470610.59

$ws.Range("H137").Value = 54346520
$ws.Range("I137").Value = 111114340
$ws.Range("J137").Value = 3255483.5
$ws.Range("K137").Value = 333343020
$ws.Range("L137").Value = 9766450.5
$ws.Range("M137").Value = -333340470
$ws.Range("N137").Value = -9771550.5


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 33060.688
$ws.Range("I32").Value = 33060.688
$ws.Range("K32").Value = 33060.688
$ws.Range("M32").Value = -32773.688

$ws.Range("H61").Value = 2568068.2
$ws.Range("I61").Value = 2781657.2
$ws.Range("K61").Value = 2781657.2
$ws.Range("M61").Value = -2781445.2

$ws.Range("H74").Value = 4352.357
$ws.Range("I74").Value = 2693.1667
$ws.Range("J74").Value = 5596.75
$ws.Range("K74").Value = 2693.1667
$ws.Range("L74").Value = 5596.75
$ws.Range("M74").Value = -1819.1667
$ws.Range("N74").Value = -7344.75

$ws.Range("H77").Value = 4352.357
$ws.Range("I77").Value = 2693.1667
$ws.Range("J77").Value = 5596.75
$ws.Range("K77").Value = 13465.8335
$ws.Range("L77").Value = 27983.75
$ws.Range("M77").Value = -9097.833500000001
$ws.Range("N77").Value = -36719.75

$ws.Range("H97").Value = 1425.9474
$ws.Range("I97").Value = 999
$ws.Range("J97").Value = 2013
$ws.Range("K97").Value = 999
$ws.Range("L97").Value = 2013
$ws.Range("M97").Value = -503
$ws.Range("N97").Value = -3005

$ws.Range("H113").Value = 99431.836
$ws.Range("J113").Value = 99431.836
$ws.Range("L113").Value = 99431.836
$ws.Range("N113").Value = -108109.836

$ws.Range("H122").Value = 2104.7144
$ws.Range("I122").Value = 2104.7144
$ws.Range("K122").Value = 6314.1432
$ws.Range("M122").Value = -3864.1432

$ws.Range("H132").Value = 771703.0600000001
$ws.Range("I132").Value = 1001814
$ws.Range("J132").Value = 4666.5
$ws.Range("K132").Value = 3005442
$ws.Range("L132").Value = 13999.5
$ws.Range("M132").Value = -3002912
$ws.Range("N132").Value = -19059.5

$ws.Range("H136").Value = 2568068.2
$ws.Range("I136").Value = 2781657.2
$ws.Range("K136").Value = 8344971.600000001
$ws.Range("M136").Value = -8342421.600000001


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1833.3334
$ws.Range("J94").Value = 1600
$ws.Range("L94").Value = 1600
$ws.Range("N94").Value = -2502

$ws.Range("H134").Value = 1015466.7
$ws.Range("I134").Value = 954088.4399999999
$ws.Range("K134").Value = 2862265.32
$ws.Range("M134").Value = -2859730.32


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 36097.617
$ws.Range("I31").Value = 13994.286
$ws.Range("J31").Value = 61884.832
$ws.Range("K31").Value = 13994.286
$ws.Range("L31").Value = 61884.832
$ws.Range("M31").Value = -13699.286
$ws.Range("N31").Value = -62474.832

$ws.Range("H34").Value = 36097.617
$ws.Range("I34").Value = 13994.286
$ws.Range("J34").Value = 61884.832
$ws.Range("K34").Value = 13994.286
$ws.Range("L34").Value = 61884.832
$ws.Range("M34").Value = -13792.286
$ws.Range("N34").Value = -62288.832

$ws.Range("H58").Value = 775778.5
$ws.Range("I58").Value = 1374047.8
$ws.Range("J58").Value = 6575.143
$ws.Range("K58").Value = 1374047.8
$ws.Range("L58").Value = 6575.143
$ws.Range("M58").Value = -1373844.8
$ws.Range("N58").Value = -6981.143

$ws.Range("H132").Value = 59127988
$ws.Range("I132").Value = 83336070
$ws.Range("J132").Value = 1028584.6
$ws.Range("K132").Value = 250008210
$ws.Range("L132").Value = 3085753.8
$ws.Range("M132").Value = -250005680
$ws.Range("N132").Value = -3090813.8

$ws.Range("H134").Value = 10767.533
$ws.Range("I134").Value = 11904.846
$ws.Range("K134").Value = 35714.538
$ws.Range("M134").Value = -33179.538

$ws.Range("H136").Value = 775778.5
$ws.Range("I136").Value = 1374047.8
$ws.Range("J136").Value = 6575.143
$ws.Range("K136").Value = 4122143.4
$ws.Range("L136").Value = 19725.429
$ws.Range("M136").Value = -4119593.4


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 8867
$ws.Range("I3").Value = 8442.25
$ws.Range("K3").Value = 25326.75
$ws.Range("M3").Value = -25214.75

$ws.Range("H86").Value = 935.9167
$ws.Range("I86").Value = 998.2727
$ws.Range("J86").Value = 250
$ws.Range("K86").Value = 2994.8181
$ws.Range("L86").Value = 750
$ws.Range("M86").Value = -1808.8181
$ws.Range("N86").Value = -3122

$ws.Range("H89").Value = 935.9167
$ws.Range("I89").Value = 998.2727
$ws.Range("J89").Value = 250
$ws.Range("K89").Value = 8984.454299999999
$ws.Range("L89").Value = 2250
$ws.Range("M89").Value = -3056.454299999999
$ws.Range("N89").Value = -14106

$ws.Range("H113").Value = 777.3684
$ws.Range("I113").Value = 283.33334
$ws.Range("J113").Value = 819.7143
$ws.Range("K113").Value = 850.0000200000001
$ws.Range("L113").Value = 2459.1429
$ws.Range("M113").Value = 1319.99998
$ws.Range("N113").Value = -6799.1429

$ws.Range("H129").Value = 2275.6365
$ws.Range("J129").Value = 3896
$ws.Range("L129").Value = 11688
$ws.Range("N129").Value = -21688


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3112.625
$ws.Range("I97").Value = 3112.625
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 3112.625
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -2616.625


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").ClearContents()
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = 0

$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").ClearContents()
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = 0

$ws.Range("H132").Value = 971639.4399999999
$ws.Range("I132").Value = 1453126.1
$ws.Range("K132").Value = 4359378.300000001
$ws.Range("M132").Value = -4356848.300000001


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 8725635
$ws.Range("I132").Value = 9789210
$ws.Range("K132").Value = 29367630
$ws.Range("M132").Value = -29365100

$ws.Range("H136").Value = 9812019
$ws.Range("I136").Value = 11913211
$ws.Range("J136").Value = 6453.8335
$ws.Range("K136").Value = 35739633
$ws.Range("L136").Value = 19361.5005
$ws.Range("M136").Value = -35737083
$ws.Range("N136").Value = -24461.5005

